$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extra "Personal" figures (rows 28-30, column M) ---
$ws.Range("M28").Value = 4648721907
$ws.Range("M29").Value = 1578318
$ws.Range("M30").Value = 1549573969
$ws.Range("M28:M30").NumberFormat = "#,##0"

# --- New data block used by the second chart (rows 32-35) ---
$ws.Range("E32").Value = "Algoritmo"
$ws.Range("F32").Value = "Mil"
$ws.Range("G32").Value = "Cinco mil"
$ws.Range("H32").Value = "Diez mil"
$ws.Range("I32").Value = "25 mil"
$ws.Range("J32").Value = "50 mil"
$ws.Range("K32").Value = "75 mil"
$ws.Range("L32").Value = "Cien mil"

$ws.Range("E33").Value = "Bubble"
$bubble = @(43751, 39154, 1389382, 9890081, 40443317, 88981014, 156658136)
for ($i = 0; $i -lt $bubble.Length; $i++) {
    $ws.Cells.Item(33, 6 + $i).Value = $bubble[$i]
}

$ws.Range("E34").Value = "Quick"
$quick = @(45746, 16805, 11084, 28931, 8129, 9282, 122297)
for ($i = 0; $i -lt $quick.Length; $i++) {
    $ws.Cells.Item(34, 6 + $i).Value = $quick[$i]
}

$ws.Range("E35").Value = "Merge"
$merge = @(18728, 184999, 645114, 4240071, 5569317, 20179699, 27681576)
for ($i = 0; $i -lt $merge.Length; $i++) {
    $ws.Cells.Item(35, 6 + $i).Value = $merge[$i]
}

$ws.Range("F33:L35").NumberFormat = "#,##0"

# --- Column widths for the new columns L (12) and M (13) ---
$ws.Columns.Item(12).ColumnWidth = 12.85546875
$ws.Columns.Item(13).ColumnWidth = 12.7109375

# --- Chart1 (existing "Bubblesort/Quicksort/Mergesort" chart): add a title ---
$chart1 = $ws.ChartObjects().Item(1).Chart
$chart1.HasTitle = $true
$chart1.ChartTitle.Text = "Personal"

# --- Chart2 (new line chart for the Bubble/Quick/Merge block) ---
$co2 = $ws.ChartObjects().Add(600, 300, 300, 260)
$co2.Name = "Grafico 1"
$chart2 = $co2.Chart
$chart2.ChartType = 4
$chart2.SetSourceData($ws.Range("E32:M35"))
$chart2.HasTitle = $true
$chart2.ChartTitle.Text = "Umd pc"

Write-Host "done"
